# Update column G ("K") values in Sheet1 per regenerated save_data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newValues = @{
    2  = 1
    3  = 1
    4  = 0
    5  = 2
    6  = 1
    7  = 0
    8  = 0
    9  = 1
    10 = 0
    11 = 0
    12 = 1
    13 = 1
    14 = 1
    15 = 0
    16 = 0
    17 = 2
    19 = 1
    20 = 1
    21 = 1
    22 = 0
    23 = 1
    24 = 1
    25 = 2
    26 = 1
    27 = 3
    28 = 2
    29 = 2
    30 = 2
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $newValues[$row]
}
